# Pushing to the new repo: add a "Sheet2" worksheet with a search-item list,
# make it the active sheet/tab, and keep Sheet1 as the first (now non-selected) tab.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# --- Add Sheet2 right after Sheet1 -----------------------------------------
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# --- Values -------------------------------------------------------------
$ws2.Range("A1").Value = "searchItem"
$ws2.Range("A2").Value = "Pencil"
$ws2.Range("A3").Value = "Iphone"
$ws2.Range("A4").Value = "Toys"
$ws2.Range("A5").Value = "women's clothing"
$ws2.Range("A6").Value = "shoes"

# --- Column width ---------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 11

# --- Formatting ------------------------------------------------------------
# A1 ("searchItem") gets the same look as the Sheet1 header row: bold font,
# thin border all around, and wrapped text. Re-use that look by copying the
# format from Sheet1's header cell (keeps the bordered look with its real
# border color instead of a brand new colorless border).
$sheet1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

# A2 ("Pencil") starts from the same bordered look, then drop the bold font
# and the wrapping so it ends up as a plain bordered cell.
$sheet1.Range("A1").Copy()
$ws2.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("A2").Font.Bold = $false
$ws2.Range("A2").WrapText = $false

# A3, A4, A6 ("Iphone", "Toys", "shoes") share that exact same plain bordered
# look as A2.
$ws2.Range("A2").Copy()
$ws2.Range("A3:A4").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("A6").PasteSpecial(-4122)     # xlPasteFormats

# A5 ("women's clothing") keeps the border and the wrapped text, but not the
# bold font.
$sheet1.Range("A1").Copy()
$ws2.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("A5").Font.Bold = $false

# Row heights for the two wrapped rows.
$ws2.Range("A1").RowHeight = 28.8
$ws2.Range("A5").RowHeight = 28.8

$ws2.Range("A1").Select()

# --- Sheet1 is no longer the displayed/active tab; Sheet2 is --------------
$ws2.Activate()

Write-Host "done"
